$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.858.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.335.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.56%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.332.24"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.25%  "
$ws.Range("E10").Value = "  +7.12%  "
$ws.Range("E11").Value = "  +2.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.73%  "
$ws.Range("E13").Value = "  +3.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "695.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.875.82"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.814.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.333.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.31%  "
$ws.Range("E22").Value = "  +1.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.62%  "
$ws.Range("E27").Value = "  +3.02%  "
$ws.Range("E28").Value = "  +6.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "566.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.106"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.693.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.134"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.56%  "
$ws.Range("E41").Value = "  +3.94%  "
$ws.Range("E42").Value = "  +8.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0676"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.89%  "
$ws.Range("E44").Value = "  +4.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0409"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.60%  "
$ws.Range("E47").Value = "  +6.12%  "
$ws.Range("E48").Value = "  +2.12%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "131.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.53%  "
